# Add exporting operations for the period
# Remove rows that fall outside the desired export period, shifting the
# remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows (bottom-up so row numbers stay valid while deleting):
# row 10 -> 27.08.2015 21:09 / category 8
# row 8  -> 11.11.2020 18:15 / category 6
# row 7  -> 15.12.2019 21:09 / category 1
# row 3  -> 29.12.2020 21:09 / category 2
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(3).Delete()
